{"js": "// Fix the trainer heading: \"About the Trainer \u2013 Mr. Krishna (Owner: Maruthi Rao)\"\n// becomes \"About the Trainer, lecturer or teacher\u2013 Mr. Krishna (Owner: Maruthi Rao)\".\n// (commit: \"prompt issue with teacher name\")\n//\n// We locate the unique phrase \"Trainer \u2013 Mr.\" (the heading's original wording,\n// using the literal en dash U+2013) and replace it in place so the surrounding\n// run formatting (bold) is preserved.\n\nconst searchText = \"Trainer \\u2013 Mr.\";\nconst replacementText = \"Trainer, lecturer or teacher\\u2013 Mr.\";\n\nconst results = context.document.body.search(searchText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"' + searchText + '\" in the document body.');\n}\n\nresults.items[0].insertText(replacementText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Fix the trainer heading: \"About the Trainer \u2013 Mr. Krishna (Owner: Maruthi Rao)\"\n# becomes \"About the Trainer, lecturer or teacher\u2013 Mr. Krishna (Owner: Maruthi Rao)\".\n# (commit: \"prompt issue with teacher name\")\n#\n# Locate the unique phrase \"Trainer \u2013 Mr.\" (using the literal en dash U+2013)\n# via Find and overwrite just that range's text so the surrounding run\n# formatting (bold) is preserved.\n\n$d = $word.ActiveDocument\n\n$enDash = [char]0x2013\n$searchText = \"Trainer \" + $enDash + \" Mr.\"\n$replacementText = \"Trainer, lecturer or teacher\" + $enDash + \" Mr.\"\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = $searchText\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$found = $find.Execute()\n\nif ($found) {\n    $range.Text = $replacementText\n} else {\n    throw \"Could not find '$searchText' in the document.\"\n}\n"}
